$wb = $excel.ActiveWorkbook
$logs = $wb.Worksheets.Item("Logs")
$dashboard = $wb.Worksheets.Item("Dashboard")

# --- Logs sheet: append new row 22 (Testmail #20) ---
$logs.Range("A22").Value = "Ik ben niet tevreden over hoe dit is gegaan."
$logs.Range("B22").Value = "mailmind.test@zohomail.eu"
$logs.Range("C22").Value = "Testmail #20: Ik ben niet tevreden over hoe dit is gegaan."
$logs.Range("D22").Value = "Klacht / Probleem"
$logs.Range("E22").Value = "Beste klant,
Dank u voor uw e-mail. Ik begrijp dat u niet tevreden bent over het proces. Om uw zorgen beter te kunnen begrijpen en aanpakken, zou u meer specifieke details kunnen delen over wat er precies is misgegaan?
Met vriendelijke groet,
[Naam bedrijf] E-mailassistent"
$logs.Range("F22").Value = "2025-07-27 19:54:38"
$logs.Range("G22").Value = "Ja"
$logs.Range("H22").Value = "Nee"
$logs.Range("I22").Value = "Ja"
$logs.Range("J22").Value = "Nee"

# The multi-line E22 value auto-sets an explicit row height on write;
# AutoFit() clears the explicit/custom height again so row 22 matches
# the other (auto-height) rows.
$logs.Rows.Item(22).AutoFit()

# --- Dashboard sheet: append new row 8 (Klacht / Probleem, 1) ---
$dashboard.Range("A8").Value = "Klacht / Probleem"
$dashboard.Range("B8").Value = 1

# --- Grow the conditional-formatting ranges on Logs from row 21 to row 22,
#     preserving dxfId/priority/grouping by moving the existing rules'
#     applies-to range instead of deleting and re-adding them. ---
$newRangeD = $logs.Range("D2:D22")
$condsD = $logs.Range("D2:D21").FormatConditions
for ($i = 1; $i -le $condsD.Count; $i++) {
    $condsD.Item($i).ModifyAppliesToRange($newRangeD)
}

$newRangeG = $logs.Range("G2:G22")
$condsG = $logs.Range("G2:G21").FormatConditions
for ($i = 1; $i -le $condsG.Count; $i++) {
    $condsG.Item($i).ModifyAppliesToRange($newRangeG)
}

$newRangeH = $logs.Range("H2:H22")
$condsH = $logs.Range("H2:H21").FormatConditions
for ($i = 1; $i -le $condsH.Count; $i++) {
    $condsH.Item($i).ModifyAppliesToRange($newRangeH)
}

$newRangeI = $logs.Range("I2:I22")
$condsI = $logs.Range("I2:I21").FormatConditions
for ($i = 1; $i -le $condsI.Count; $i++) {
    $condsI.Item($i).ModifyAppliesToRange($newRangeI)
}

$newRangeJ = $logs.Range("J2:J22")
$condsJ = $logs.Range("J2:J21").FormatConditions
for ($i = 1; $i -le $condsJ.Count; $i++) {
    $condsJ.Item($i).ModifyAppliesToRange($newRangeJ)
}

# --- Update the Dashboard bar chart's category/value series to extend
#     through the new row (A2:A8 / B2:B8). ---
$chartObj = $dashboard.ChartObjects(1)
$chart = $chartObj.Chart
$ser = $chart.SeriesCollection(1)
$ser.Values = "='Dashboard'!`$B`$2:`$B`$8"
$ser.XValues = "='Dashboard'!`$A`$2:`$A`$8"
